$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source cells contain a mojibake rendering of the UTF-8 "±" (U+00B1)
# sign: it was double-encoded as U+00C2 U+00B1 ("Â±"). Restore the correct
# single "±" character across the affected result columns.
$mojibake = [string]([char]0x00C2) + [string]([char]0x00B1)
$plusMinus = [string]([char]0x00B1)

for ($row = 2; $row -le 17; $row++) {
    foreach ($col in @("B", "C", "D")) {
        $cell = $ws.Range("$col$row")
        $text = $cell.Text
        if ($text -ne $null -and $text.Contains($mojibake)) {
            $cell.Value = $text.Replace($mojibake, $plusMinus)
        }
    }
}
